$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells to Text format so that numeric-looking
# strings (e.g. "1.006") are preserved exactly as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Column B (Coin name) updates ---
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("B40").Value = 'Frax'
$ws.Range("B41").Value = 'Aptos'
$ws.Range("B42").Value = 'Algorand'

# --- Column C (Link) updates ---
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'

# --- Column D (Price) updates ---
$ws.Range("D2").Value = '20.228.14'
$ws.Range("D3").Value = '1.441.93'
$ws.Range("D5").Value = '0.9196'
$ws.Range("D6").Value = '274.51'
$ws.Range("D7").Value = '0.3629'
$ws.Range("D8").Value = '0.3073'
$ws.Range("D9").Value = '38.74'
$ws.Range("D10").Value = '1.018'
$ws.Range("D11").Value = '0.06479'
$ws.Range("D12").Value = '0.9986'
$ws.Range("D13").Value = '5.322'
$ws.Range("D14").Value = '17.39'
$ws.Range("D15").Value = '6.019'
$ws.Range("D16").Value = '0.00001006'
$ws.Range("D17").Value = '1.440.16'
$ws.Range("D18").Value = '0.9360'
$ws.Range("D19").Value = '0.05605'
$ws.Range("D20").Value = '67.45'
$ws.Range("D21").Value = '5.369'
$ws.Range("D23").Value = '10.77'
$ws.Range("D24").Value = '2.240'
$ws.Range("D25").Value = '20.249.25'
$ws.Range("D26").Value = '138.36'
$ws.Range("D27").Value = '2.055'
$ws.Range("D28").Value = '16.87'
$ws.Range("D29").Value = '1.593.90'
$ws.Range("D30").Value = '109.93'
$ws.Range("D31").Value = '3.984'
$ws.Range("D32").Value = '4.818'
$ws.Range("D33").Value = '0.7855'
$ws.Range("D34").Value = '0.07629'
$ws.Range("D35").Value = '1.459'
$ws.Range("D36").Value = '0.05782'
$ws.Range("D37").Value = '1.137'
$ws.Range("D38").Value = '4.633'
$ws.Range("D39").Value = '0.01978'
$ws.Range("D40").Value = '0.9318'
$ws.Range("D41").Value = '10.11'
$ws.Range("D42").Value = '0.1840'
$ws.Range("D43").Value = '6.993'
$ws.Range("D44").Value = '0.5173'
$ws.Range("D45").Value = '3.472'
$ws.Range("D46").Value = '11.69'
$ws.Range("D47").Value = '115.68'
$ws.Range("D48").Value = '0.5078'
$ws.Range("D49").Value = '1.725'
$ws.Range("D50").Value = '0.06355'
$ws.Range("D51").Value = '0.9880'

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("E3").Value = '  +2.66%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("E5").Value = '  -8.14%  '
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("E7").Value = '  -1.33%  '
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("E13").Value = '  -1.77%  '
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("E15").Value = '  -2.15%  '
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("E17").Value = '  +2.69%  '
$ws.Range("E18").Value = '  -6.51%  '
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("E20").Value = '  -4.22%  '
$ws.Range("E21").Value = '  -3.48%  '
$ws.Range("E22").Value = '  -3.41%  '
$ws.Range("E23").Value = '  -2.14%  '
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("E27").Value = '  -8.08%  '
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("E29").Value = '  +2.12%  '
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("E31").Value = '  -2.50%  '
$ws.Range("E32").Value = '  -8.83%  '
$ws.Range("E33").Value = '  -3.04%  '
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("E35").Value = '  +1.29%  '
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("E37").Value = '  +4.83%  '
$ws.Range("E38").Value = '  -3.62%  '
$ws.Range("E39").Value = '  -4.09%  '
$ws.Range("E40").Value = '  -6.93%  '
$ws.Range("E41").Value = '  -2.30%  '
$ws.Range("E42").Value = '  -3.32%  '
$ws.Range("E43").Value = '  -16.53%  '
$ws.Range("E44").Value = '  -1.69%  '
$ws.Range("E45").Value = '  -0.84%  '
$ws.Range("E46").Value = '  -4.31%  '
$ws.Range("E47").Value = '  +3.46%  '
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("E50").Value = '  +3.22%  '
$ws.Range("E51").Value = '  -1.40%  '
